$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add the two new data rows (Y Bot Model / Male Locomotion Animation Pack)
#    Values are entered in a specific column order so that the resulting
#    shared-string table is built up in the same order as the target file.
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Male Locomotion Animation Pack"
$ws.Range("A5").Value = "Y Bot Model"
$ws.Range("B5").Value = "Y Bot Model"
$ws.Range("B6").Value = "Male Locomotion Animation Pack"
$ws.Range("C5").Value = "Mixamo"
$ws.Range("C6").Value = "Mixamo"
$ws.Range("D5").Value = "Adobe Terms of Use (Free for video game use)"
$ws.Range("D6").Value = "Adobe Terms of Use (Free for video game use)"
$ws.Range("E5").Value = "3D Model(s)"
$ws.Range("E6").Value = "Animation(s)"
$ws.Range("F5").Value = "3D Model we are using for our player character"
$ws.Range("F6").Value = "Animations for our player character"

# ---------------------------------------------------------------------------
# 2. Apply word-wrap formatting to the whole table area (header + all data
#    rows) so that every column ends up using the same centered / wrapped
#    alignment that column F already had.
# ---------------------------------------------------------------------------
$ws.Range("A1:F6").WrapText = $true
$ws.Range("A1:F6").HorizontalAlignment = -4108
$ws.Range("A1:F6").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 3. Row heights - rows with wrapped, two-line text grow to 30pt.
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(6).RowHeight = 30

# ---------------------------------------------------------------------------
# 4. Column widths - nudged slightly, matching the after-edit workbook.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 23.57
$ws.Columns.Item(2).ColumnWidth = 22.71
$ws.Columns.Item(3).ColumnWidth = 23.57
$ws.Columns.Item(4).ColumnWidth = 30.86
$ws.Columns.Item(5).ColumnWidth = 12.86
$ws.Columns.Item(6).ColumnWidth = 35.86

# ---------------------------------------------------------------------------
# 5. Selection / active cell, matching the saved view state.
# ---------------------------------------------------------------------------
[void]$ws.Range("D10").Select()

Write-Host "Edit complete"
